$p = $ppt.ActivePresentation

# Slide 1: notes page text "Some speaker notes" (merge runs into one)
$s1 = $p.Slides.Item(1)
$notes1 = $s1.NotesPage
$notes1.Shapes.Item(2).TextFrame.TextRange.Text = "Some speaker notes"

# Slide 1: subtitle text "Jesse Rosenthal" (merge runs into one, keep leading breaks)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "`r`rJesse Rosenthal"

# Slide 2: title text "A header" (merge runs into one)
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "A header"
